$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    $ws.Range("D$r").Value = 0.169
    $ws.Range("E$r").Value = 0.6679999999999999
    $ws.Range("F$r").Value = 0.0239
    $ws.Range("G$r").Value = 0.2109645265318089
    $ws.Range("H$r").Value = 0.2109645265318089
    $ws.Range("I$r").Value = 0.236880680152448
    $ws.Range("J$r").Value = 0.1725036975548375
    $ws.Range("K$r").Value = 310.9
    $ws.Range("L$r").Value = 0.1822925828202873
    $ws.Range("M$r").Value = 64.90000000000001
    $ws.Range("N$r").Value = 0.02962793882675189
    $ws.Range("O$r").Value = 0.2087487938243809
    $ws.Range("P$r").Value = 32.7
    $ws.Range("Q$r").Value = 0.01492809860762383
    $ws.Range("R$r").Value = 0.1051785139916372
    $ws.Range("S$r").Value = 32.2
    $ws.Range("T$r").Value = 0.4961479198767335
    $ws.Range("U$r").Value = 72
    $ws.Range("V$r").Value = 0.03286920794339192
    $ws.Range("W$r").Value = 0.5028303412582888
    $ws.Range("X$r").Value = 0.05032700534849797
    $ws.Range("Y$r").Value = 0.4525033359097909
    $ws.Range("Z$r").Value = 3.378565768621236
    $ws.Range("AA$r").Value = 0.5828150875193647
    $ws.Range("AB$r").Value = 0.05032700534849797
    $ws.Range("AC$r").Value = 0.5324880821708667
    $ws.Range("AG$r").Value = -72
    $ws.Range("AJ$r").Value = -0.0339863110691527
    $ws.Range("AK$r").Value = -0.09770660876645407
    $ws.Range("AP$r").Value = -0.1691729323308271
}
